$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Four new "measurement" blocks appended below the existing data (rows 19-22).
# Column A holds the SQL query used for the measurement, columns B-I hold the
# individual timing results. The timing values are numeric-looking strings
# that must be stored as text (as in the rest of the sheet), so each such
# cell is briefly switched to a text number format while its value is set and
# then restored to the Normal style so it keeps the sheet's default look.

$rows = @(
    @("select id from product where retailprice_cs = 110000",
      "0.00968390","0.00091721","0.00074059","0.00047906",
      "0.00070678","0.00050203","0.00067805","0.00053893"),
    @("select id from product where retailprice_cs = 2932650",
      "0.00050963","0.00051144","0.00767151","0.00047020",
      "0.00053813","0.00068424","0.00049065","0.00051151"),
    @("select id from product where retailprice_cs = 1855875",
      "0.00050368","0.00046416","0.00046276","0.00049457",
      "0.00046781","0.00049790","0.00059695","0.00063842"),
    @("select id from product where retailprice_cs = 954375",
      "0.00062934","0.00064824","0.00058874","0.00057954",
      "0.00066041","0.00068109","0.00054907","0.00071128")
)

$columns = @("A","B","C","D","E","F","G","H","I")
$startRow = 19

for ($i = 0; $i -lt $rows.Count; $i++) {
    $rowIndex = $startRow + $i
    $values = $rows[$i]
    for ($j = 0; $j -lt $values.Count; $j++) {
        $cell = $ws.Range($columns[$j] + $rowIndex)
        if ($j -eq 0) {
            $cell.Value = $values[$j]
        } else {
            $cell.NumberFormat = "@"
            $cell.Value = $values[$j]
            $cell.Style = "Normal"
        }
    }
}

$ws.Range("I23").Select() | Out-Null
